# Auto-optimize exam scheduling: dynamically adjusts exams per slot (1-4)
# to guarantee all courses are scheduled within date range.
# This updates the Section_A / Section_B weekly timetables and the
# Elective_Coordination schedule for CS151 accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Section_A
# ---------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Section_A")

$wsA.Range("B2").Value = "CS151 (Elective)"
$wsA.Range("C2").Value = "MA162"
$wsA.Range("D2").Value = "Free"
$wsA.Range("E2").Value = "CS161"
$wsA.Range("F2").Value = "CS151 (Elective)"

$wsA.Range("B3").Value = "MA162"
$wsA.Range("D3").Value = "CS161"
$wsA.Range("E3").Value = "Free"

$wsA.Range("B5").Value = "MA161"
$wsA.Range("C5").Value = "EC161"
$wsA.Range("F5").Value = "CS161"

$wsA.Range("B6").Value = "CS151 (Tutorial)"

$wsA.Range("C7").Value = "CS161"
$wsA.Range("D7").Value = "Free"
$wsA.Range("E7").Value = "C202"
$wsA.Range("F7").Value = "EC161"

$wsA.Range("E8").Value = "Free"

# ---------------------------------------------------------------------
# Section_B
# ---------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("B2").Value = "CS151 (Elective)"
$wsB.Range("C2").Value = "Free"
$wsB.Range("E2").Value = "MA162"
$wsB.Range("F2").Value = "CS151 (Elective)"

$wsB.Range("B3").Value = "MA161"
$wsB.Range("C3").Value = "Free"
$wsB.Range("D3").Value = "Free"
$wsB.Range("E3").Value = "C202"
$wsB.Range("F3").Value = "EC161"

$wsB.Range("B5").Value = "CS161"
$wsB.Range("C5").Value = "C202"
$wsB.Range("D5").Value = "MA162"
$wsB.Range("E5").Value = "MA161"

$wsB.Range("B6").Value = "CS151 (Tutorial)"

$wsB.Range("C7").Value = "CS161"
$wsB.Range("D7").Value = "CS161"
$wsB.Range("E7").Value = "EC161"
$wsB.Range("F7").Value = "C202"

$wsB.Range("E8").Value = "Free"

# ---------------------------------------------------------------------
# Elective_Coordination (CS151 lecture/tutorial slots)
# ---------------------------------------------------------------------
$wsE = $wb.Worksheets.Item("Elective_Coordination")

$wsE.Range("C10").Value = "Fri"

$wsE.Range("C11").Value = "Mon"
$wsE.Range("D11").Value = "09:00-10:30"

$wsE.Range("C12").Value = "Mon"
$wsE.Range("D12").Value = "14:30-15:30"
